$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Original layout:
#   1 Vingadores
#   2 Superman
#   3 Batman
#   4 "Starwars" (wrapped in proofErr spellcheck tags) + a separate " " run
#   5 "Harry " run + "potter" run (wrapped in proofErr spellcheck tags)
#   6 Velozes e furiosos
#   7 (empty)
#   8 (empty)
#
# Target layout:
#   1 Vingadores
#   2 Superman
#   3 Batman
#   4 "Starwars " (single clean run, no spellcheck markup)
#   5 "Harry potter" (new paragraph, single clean run)
#   6 Velozes e furiosos + a trailing " " run
#   7 "Planeta dos macacos" (new paragraph)
#   8 (empty)
#   9 (empty)
# ---------------------------------------------------------------------------

# --- Paragraph 4: "Starwars" -> clean single-run "Starwars " paragraph -----
# Deleting the paragraph's range (not the mark) removes its runs/proofErr
# markers and merges the (now empty) paragraph with the next one; we then
# retype clean text followed by a paragraph break to recreate the paragraph.
$p4 = $d.Paragraphs.Item(4)
$r4 = $p4.Range
$r4.Delete()
$ins4 = $d.Range($r4.Start, $r4.Start)
$ins4.InsertBefore("Starwars `r")

# --- Paragraph 5: "Harry " + "potter" -> clean single-run "Harry potter" --
$p5 = $d.Paragraphs.Item(5)
$r5 = $p5.Range
$r5.Delete()
$ins5 = $d.Range($r5.Start, $r5.Start)
$ins5.InsertBefore("Harry potter`r")

# --- Paragraph 6: "Velozes e furiosos" gains a trailing separate " " run -
# Re-use the genuinely separate, un-formatted trailing-space run that
# already lives after "Starwars " (paragraph 4) so the new run is a real
# distinct run (matching the source document's own pattern) instead of
# being silently merged into the preceding text run.
$p4b = $d.Paragraphs.Item(4)
$full4 = $p4b.Range
$spaceChar = $full4.Characters.Item($full4.Characters.Count - 1)
$ft = $spaceChar.FormattedText

$p6 = $d.Paragraphs.Item(6)
$r6 = $p6.Range
$r6.Collapse(0)
$r6.FormattedText = $ft

# --- New paragraph after "Velozes e furiosos": "Planeta dos macacos" -----
$p6b = $d.Paragraphs.Item(6)
$p6b.Range.InsertParagraphAfter()
$p7 = $d.Paragraphs.Item(7)
$p7.Range.InsertBefore("Planeta dos macacos")
